# "Problem statement and features" commit:
#   - Heading "Features: -" becomes "Features/Business Cases: -"
#   - The five existing bullet items are rewritten/reshuffled and a new
#     sixth bullet ("Facility for cancellation of tickets.") is appended.
#   - The hidden "_GoBack" bookmark (Word's "last edit location" marker)
#     ends up inside the "Searching for Trains ..." bullet instead of at
#     the end of the list.
#
# Each bullet is addressed by its paragraph index so the text shuffle
# (values moving from one bullet position to another) can't collide with
# itself the way a sequence of Find/Replace-on-text calls could.

$d = $word.ActiveDocument

# --- Heading: "Features: -" -> "Features/Business Cases: -" ----------
$p = $d.Paragraphs.Item(15)
$r = $p.Range
[void]$r.MoveEnd(1, -1)
$r.Text = "Features/Business Cases: -"

# --- Bullet 1: add the new lead-in sentence fragment ------------------
$p = $d.Paragraphs.Item(16)
$r = $p.Range
[void]$r.MoveEnd(1, -1)
$r.Text = "The system should be able to add Train Routes consisting of all the middle stations."

# --- Bullet 2: brand-new login/roles feature ---------------------------
$p = $d.Paragraphs.Item(17)
$r = $p.Range
[void]$r.MoveEnd(1, -1)
$r.Text = "System should allow customers and Train officials to login and perform various operations."

# --- Bullet 3: old bullet-2 text moves here ----------------------------
$p = $d.Paragraphs.Item(18)
$r = $p.Range
[void]$r.MoveEnd(1, -1)
$r.Text = "Searching for Trains from source to destination on particular date."

# --- Bullet 4: old bullet-3 text moves here ----------------------------
$p = $d.Paragraphs.Item(19)
$r = $p.Range
[void]$r.MoveEnd(1, -1)
$r.Text = "Checking for availability of seats in the train."

# --- Bullet 5: old bullet-4 text moves here (bookmark no longer here) --
$p = $d.Paragraphs.Item(20)
$r = $p.Range
[void]$r.MoveEnd(1, -1)
$r.Text = "View/ Print booked Train Tickets. Also view past bookings."

# Drop the old "_GoBack" bookmark that used to sit at the end of the
# (now retired) "Facility for cancellation" paragraph; it gets re-added
# further down, inside bullet 3.
$bm = $d.Bookmarks.Item("_GoBack")
$bm.Delete()

# --- Bullet 6 (new): old bullet-5 text re-appears as its own paragraph -
$p20 = $d.Paragraphs.Item(20)
$p20.Range.InsertParagraphAfter()
$p21 = $d.Paragraphs.Item(21)
$r21 = $p21.Range
[void]$r21.MoveEnd(1, -1)
$r21.Text = "Facility for cancellation of tickets."

# Re-anchor "_GoBack" between "Se" and "arching" inside bullet 3's text.
$p18 = $d.Paragraphs.Item(18)
$bmPos = $p18.Range.Start + 2
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)
